# Extend the 2x2 "gradient descent method" grid display from columns A:B
# out to columns A:E, mirroring the existing layout (value 0, same column
# width) into the three newly-added columns C, D and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values (rows 1-2, columns C-E), matching the existing A:B cells.
$ws.Range("C1:E2").Value = 0

# Match the narrow column width already used for columns A and B
# (~2.14 "characters" wide).
$ws.Columns.Item(3).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(4).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(5).ColumnWidth = 1.3333333333333333
